$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Notes" worksheet right after "Dashboard"
# ---------------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$notes = $wb.Worksheets.Add($null, $dashboard)
$notes.Name = "Notes"

# Column width for column A
$notes.Columns.Item(1).ColumnWidth = 38

# ---- values / formulas first (styling is applied afterwards so it cannot
#      "bleed" onto cells whose formulas reference an already-styled cell) --

# Row 2
$notes.Range("A2").Value = "Electricity demand services "

# Row 3 - headers
$notes.Range("B3").Value = "Unit"
$notes.Range("D3").Value = "Value "
$notes.Range("F3").Value = "Source"

# Row 4
$notes.Range("A4").Value = "Electricity use services Nieuw-West (wijk)"
$notes.Range("B4").Value = "kWh"
$notes.Range("D4").Value = 21530000
$notes.Range("F4").Value = "Klimaatmonitor: https://klimaatmonitor.databank.nl/Jive?workspace_guid=0260472e-5a67-45b9-9097-f3173bf02a73"

# Row 5
$notes.Range("A5").Value = "Installed capacity Solar PV "
$notes.Range("B5").Value = "kWp"
$notes.Range("D5").Formula = "=Dashboard!F18"
$notes.Range("F5").Value = "Dashboard"

# Row 6
$notes.Range("A6").Value = "Full load hours PV"
$notes.Range("B6").Value = "h"
$notes.Range("D6").Formula = "=D5*Dashboard!F111"
$notes.Range("F6").Value = "ETM"

# Row 8
$notes.Range("A8").Value = "Registered companies Reitdiep"
$notes.Range("B8").Value = "#"
$notes.Range("D8").Value = 160
$notes.Range("F8").Value = "CBS - gevestigde bedrijven"

# Row 9
$notes.Range("A9").Value = "Registered companies Nieuw-West"
$notes.Range("B9").Value = "#"
$notes.Range("D9").Value = 790
$notes.Range("F9").Value = "CBS - gevestigde bedrijven"

# Row 11 - result row
$notes.Range("A11").Value = "Electricity use services Reitdiep"
$notes.Range("B11").Value = "kWh"
$notes.Range("D11").Formula = "=D4 * (D8/D9) +D6"
$notes.Range("F11").Value = "Geschaald Nieuw-West + solar PV"

# ---- styling, applied last ------------------------------------------------
# D6 and D11 are formulas that (indirectly) reference Dashboard cells whose
# own style carries a "#,##0" number format; the recalculation engine can
# otherwise let that format "leak" onto these dependent cells, so their
# style is reset to Normal right before the intentional formatting below is
# applied.
$notes.Range("D6").Style = "Normal"
$notes.Range("D11").Style = "Normal"

$notes.Range("A2").Font.Bold = $true
$notes.Range("B3").Font.Bold = $true
$notes.Range("D3").Font.Bold = $true
$notes.Range("F3").Font.Bold = $true
$notes.Range("A11:H11").Font.Bold = $true

$notes.Range("D5").NumberFormat = "#,##0"
$notes.Range("D8").NumberFormat = "#,##0"
$notes.Range("D9").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# 2. Dashboard!F22 now pulls its value from the Notes sheet instead of being
#    computed locally, and the helper value in G22 is cleared.
# ---------------------------------------------------------------------------
$dashboard.Range("F22").Formula = "=Notes!D11"
$dashboard.Range("G22").ClearContents()

# ---------------------------------------------------------------------------
# 3. Restore selections / active sheet state to match the saved workbook
# ---------------------------------------------------------------------------
$dashboard.Range("H10").Select()
$notes.Range("I18").Select()
$wb.Worksheets.Item("Analyse").Range("E42").Select()
$wb.Worksheets.Item(".yml").Range("A40").Select()
